$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Date: 2025-10-03T16:37:46+01:00 -> 2025-11-27T11:57:11+00:00
$meta.Range("B8").Value = "2025-11-27T11:57:11+00:00"

# Context: element:Element -> element:Observation
$meta.Range("B26").Value = "element:Observation"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Binding Value Set URL update
$elem.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/environmental-context-vs"

# Column Z width update (59.50390625 -> 62.04296875)
# The host engine snaps ColumnWidth to a pixel grid (~1/6 increments) when
# it re-serializes column metadata, so feed it the input that lands on the
# grid point closest to the target stored width (62.0 vs 62.04296875).
$elem.Columns.Item(26).ColumnWidth = 61.15
